$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 26 (pushing existing rows 26..41 down to 28..43),
# copying formatting (date style, etc.) from the row above.
$ws.Rows.Item(26).Resize(2).Insert()

# New row 26: "Primera" quality entry for 2022-06-17 (serial 44729)
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C26").Value = "Arica y Parinacota"
$ws.Range("D26").Value = 44729
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100108
$ws.Range("H26").Value = "Tropicales y subtropicales"
$ws.Range("I26").Value = 100108001
$ws.Range("J26").Value = "Guayaba"
$ws.Range("K26").Value = "Sin especificar"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 150
$ws.Range("N26").Value = 700
$ws.Range("O26").Value = 800
$ws.Range("P26").Value = 750
$ws.Range("Q26").Value = "$/kilo (en caja de 10 kilos )"
$ws.Range("R26").Value = "Región de Arica y Parinacota"
$ws.Range("S26").Value = 750
$ws.Range("T26").Value = 1

# New row 27: "Segunda" quality entry for 2022-06-17 (serial 44729)
$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").Value = 44729
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100108
$ws.Range("H27").Value = "Tropicales y subtropicales"
$ws.Range("I27").Value = 100108001
$ws.Range("J27").Value = "Guayaba"
$ws.Range("K27").Value = "Sin especificar"
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 160
$ws.Range("N27").Value = 500
$ws.Range("O27").Value = 600
$ws.Range("P27").Value = 550
$ws.Range("Q27").Value = "$/kilo (en caja de 10 kilos )"
$ws.Range("R27").Value = "Región de Arica y Parinacota"
$ws.Range("S27").Value = 550
$ws.Range("T27").Value = 1
